# إضافة حدث جديد في Card19 by HOSSAM at 2025-12-08 11:40:25
# Append a new service-log row (row 20) to the "Card19" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

$row = 20

# Column A holds the card number as text ("19"), same as every other row.
$ws.Cells.Item($row, 1).Value = "'19"

# Columns B..K (Min/Max tones + the checkmark columns) stay blank for this entry.
for ($col = 2; $col -le 11; $col++) {
    $ws.Cells.Item($row, $col).Value = ""
}

# Date of the new event. A leading apostrophe keeps it a literal text string
# (matching the other Date-column entries, which are stored as text) instead
# of letting Excel reinterpret "5/3/2025" as a date serial.
$ws.Cells.Item($row, 12).Value = "'5/3/2025"

# Event / Correction / Serviced by
$ws.Cells.Item($row, 13).Value = "قطع سير كويلر مسنن 1270"
$ws.Cells.Item($row, 14).Value = "تم تغير سير 1270"
$ws.Cells.Item($row, 15).Value = "فني"
